$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($wsArg, $addr, $value)
    $cell = $wsArg.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '22.266.65'
Set-TextValue $ws 'E2' '  -1.25%  '

Set-TextValue $ws 'D3' '1.555.84'
Set-TextValue $ws 'E3' '  -1.36%  '

Set-TextValue $ws 'D4' '1.002'
Set-TextValue $ws 'E4' '  -0.07%  '

Set-TextValue $ws 'D6' '287.54'
Set-TextValue $ws 'E6' '  -0.38%  '

Set-TextValue $ws 'D7' '0.3770'
Set-TextValue $ws 'E7' '  +2.24%  '

Set-TextValue $ws 'D8' '0.3259'
Set-TextValue $ws 'E8' '  -2.54%  '

Set-TextValue $ws 'D9' '43.93'
Set-TextValue $ws 'E9' '  -9.56%  '

Set-TextValue $ws 'D10' '1.133'
Set-TextValue $ws 'E10' '  -0.97%  '

Set-TextValue $ws 'D11' '0.07356'
Set-TextValue $ws 'E11' '  -1.56%  '

Set-TextValue $ws 'D12' '1.002'
Set-TextValue $ws 'E12' '  -0.06%  '

Set-TextValue $ws 'D13' '20.20'
Set-TextValue $ws 'E13' '  -3.81%  '

Set-TextValue $ws 'D14' '5.833'
Set-TextValue $ws 'E14' '  -2.80%  '

Set-TextValue $ws 'D15' '6.740'
Set-TextValue $ws 'E15' '  -3.18%  '

Set-TextValue $ws 'D16' '1.547.33'
Set-TextValue $ws 'E16' '  -2.11%  '

Set-TextValue $ws 'E17' '  -3.91%  '

Set-TextValue $ws 'D18' '0.06638'
Set-TextValue $ws 'E18' '  -1.86%  '

Set-TextValue $ws 'E19' '  -3.19%  '

Set-TextValue $ws 'E20' '  -0.68%  '

Set-TextValue $ws 'E21' '  -0.04%  '

Set-TextValue $ws 'D22' '16.07'
Set-TextValue $ws 'E22' '  -3.06%  '

Set-TextValue $ws 'D23' '11.63'
Set-TextValue $ws 'E23' '  -4.53%  '

Set-TextValue $ws 'D24' '22.252.59'
Set-TextValue $ws 'E24' '  -1.33%  '

Set-TextValue $ws 'D25' '2.296'
Set-TextValue $ws 'E25' '  -4.27%  '

Set-TextValue $ws 'D26' '2.553'
Set-TextValue $ws 'E26' '  -2.22%  '

Set-TextValue $ws 'D27' '151.14'
Set-TextValue $ws 'E27' '  -0.88%  '

Set-TextValue $ws 'D28' '19.31'
Set-TextValue $ws 'E28' '  -1.98%  '

Set-TextValue $ws 'D29' '4.925'
Set-TextValue $ws 'E29' '  -1.96%  '

Set-TextValue $ws 'D30' '122.45'
Set-TextValue $ws 'E30' '  -1.64%  '

Set-TextValue $ws 'D31' '1.723.50'
Set-TextValue $ws 'E31' '  -1.93%  '

Set-TextValue $ws 'D32' '1.077'
Set-TextValue $ws 'E32' '  +0.44%  '

Set-TextValue $ws 'D33' '5.893'
Set-TextValue $ws 'E33' '  -4.89%  '

Set-TextValue $ws 'D34' '1.899'
Set-TextValue $ws 'E34' '  -5.32%  '

Set-TextValue $ws 'D35' '9.366'
Set-TextValue $ws 'E35' '  -3.42%  '

Set-TextValue $ws 'D36' '0.08197'
Set-TextValue $ws 'E36' '  -1.37%  '

Set-TextValue $ws 'D37' '0.02359'
Set-TextValue $ws 'E37' '  -3.98%  '

Set-TextValue $ws 'D38' '0.06277'
Set-TextValue $ws 'E38' '  -1.56%  '

Set-TextValue $ws 'D39' '5.283'
Set-TextValue $ws 'E39' '  -3.07%  '

Set-TextValue $ws 'D40' '0.2146'
Set-TextValue $ws 'E40' '  -5.34%  '

Set-TextValue $ws 'D41' '1.246'
Set-TextValue $ws 'E41' '  -4.55%  '

Set-TextValue $ws 'E42' '  -3.42%  '

Set-TextValue $ws 'D43' '0.6054'
Set-TextValue $ws 'E43' '  -4.68%  '

Set-TextValue $ws 'D44' '1.001'
Set-TextValue $ws 'E44' '  -0.03%  '

Set-TextValue $ws 'D45' '13.74'
Set-TextValue $ws 'E45' '  -1.98%  '

Set-TextValue $ws 'D46' '3.741'
Set-TextValue $ws 'E46' '  -0.79%  '

Set-TextValue $ws 'D47' '0.5900'
Set-TextValue $ws 'E47' '  -4.52%  '

Set-TextValue $ws 'B48' 'NEARProtocol'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D48' '1.976'
Set-TextValue $ws 'E48' '  -4.18%  '

Set-TextValue $ws 'B49' 'Quant'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D49' '122.43'
Set-TextValue $ws 'E49' '  -2.15%  '

Set-TextValue $ws 'D50' '1.173'
Set-TextValue $ws 'E50' '  -3.75%  '

Set-TextValue $ws 'D51' '0.07066'
Set-TextValue $ws 'E51' '  -2.91%  '
